$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set H3 to the new "WhiteScreenBehavior" text (adds a new shared string entry).
$ws.Range("H3").Value = "WhiteScreenBehavior"

# Widen column D to fit the new, longer header text.
$ws.Columns.Item(4).ColumnWidth = 93.140625

# Row 11 gets an explicit custom height.
$ws.Rows.Item(11).RowHeight = 14.25

# Move the active cell selection to H3.
$ws.Range("H3").Select()
